$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The class schedule's afternoon block shifts earlier by one slot (lunch now
# starts at 12:20 instead of 13:00) and the day is extended with three more
# 50-minute slots (16:40, 17:30, 18:20) so the grid grows from A1:F14 to
# A1:F17. Rewrite rows 8-17 (time column + the five weekday columns) with
# their final values directly.

$ws.Range("A8").Value = "11:30"
$ws.Range("B8:F8").Value = "-"

$ws.Range("A9").Value = "12:20"
$ws.Range("B9:F9").Value = "Almoço"

$ws.Range("A10").Value = "13:00"
$ws.Range("B10:F10").Value = "-"

$ws.Range("A11").Value = "13:50"
$ws.Range("B11:F11").Value = "-"

$ws.Range("A12").Value = "14:40"
$ws.Range("B12:F12").Value = "-"

$ws.Range("A13").Value = "15:30"
$ws.Range("B13:F13").Value = "Intervalo"

$ws.Range("A14").Value = "15:50"
$ws.Range("B14:F14").Value = "-"

$ws.Range("A15").Value = "16:40"
$ws.Range("B15:F15").Value = "-"

$ws.Range("A16").Value = "17:30"
$ws.Range("B16:F16").Value = "-"

$ws.Range("A17").Value = "18:20"
$ws.Range("B17:F17").Value = ""
